$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7476426666666667
$ws.Range("H2").Value = 2.242928
$ws.Range("I2").Value = 0.001581772089386036
$ws.Range("J2").Value = 0.001581772089386036
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 2.063153084158222
$ws.Range("R2").Value = 18.568377757424
$ws.Range("S2").Value = 0.0004071587870307206
$ws.Range("T2").Value = 0.0004071587870307206
$ws.Range("G3").Value = 0.7476426666666667
$ws.Range("H3").Value = 2.242928
$ws.Range("I3").Value = 0.001581772089386036
$ws.Range("J3").Value = 0.001581772089386036
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 5.416840834885333
$ws.Range("R3").Value = 48.751567513968
$ws.Range("S3").Value = 0.001069001791871518
$ws.Range("T3").Value = 0.001069001791871518
$ws.Range("G4").Value = 0.7476426666666667
$ws.Range("H4").Value = 2.242928
$ws.Range("I4").Value = 0.001581772089386036
$ws.Range("J4").Value = 0.001581772089386036
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 0.3887976128035556
$ws.Range("R4").Value = 3.499178515232
$ws.Range("S4").Value = 0.00007672836574515445
$ws.Range("T4").Value = 0.00007672836574515445
$ws.Range("G5").Value = 0.7476426666666667
$ws.Range("H5").Value = 2.242928
$ws.Range("I5").Value = 0.001581772089386036
$ws.Range("J5").Value = 0.001581772089386036
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 0.1463565347128889
$ws.Range("R5").Value = 1.317208812416
$ws.Range("S5").Value = 0.00002888314473864289
$ws.Range("T5").Value = 0.00002888314473864289
$ws.Range("I6").Value = 0.02590993131491687
$ws.Range("J6").Value = 0.02590993131491688
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 33.79510554105644
$ws.Range("R6").Value = 304.155949869508
$ws.Range("S6").Value = 0.006669390790885433
$ws.Range("T6").Value = 0.006669390790885433
$ws.Range("I7").Value = 0.02590993131491687
$ws.Range("J7").Value = 0.02590993131491688
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("Q7").Value = 88.72958052395067
$ws.Range("R7").Value = 798.566224715556
$ws.Range("S7").Value = 0.01751059029854609
$ws.Range("T7").Value = 0.01751059029854609
$ws.Range("I8").Value = 0.02590993131491687
$ws.Range("J8").Value = 0.02590993131491688
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 6.368628900927113
$ws.Range("R8").Value = 57.31766010834401
$ws.Range("S8").Value = 0.001256835102669199
$ws.Range("T8").Value = 0.001256835102669199
$ws.Range("I9").Value = 0.02590993131491687
$ws.Range("J9").Value = 0.02590993131491688
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 2.397366717585778
$ws.Range("R9").Value = 21.576300458272
$ws.Range("S9").Value = 0.0004731151228161546
$ws.Range("T9").Value = 0.0004731151228161547
$ws.Range("G10").Value = 18.93023433333333
$ws.Range("H10").Value = 56.79070299999999
$ws.Range("I10").Value = 0.04005030430848061
$ws.Range("J10").Value = 0.04005030430848062
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 52.2388208832221
$ws.Range("R10").Value = 470.1493879489989
$ws.Range("S10").Value = 0.01030921801685203
$ws.Range("T10").Value = 0.01030921801685204
$ws.Range("G11").Value = 18.93023433333333
$ws.Range("H11").Value = 56.79070299999999
$ws.Range("I11").Value = 0.04005030430848061
$ws.Range("J11").Value = 0.04005030430848062
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 137.1538449081937
$ws.Range("R11").Value = 1234.384604173743
$ws.Range("S11").Value = 0.02706701386252398
$ws.Range("T11").Value = 0.02706701386252398
$ws.Range("G12").Value = 18.93023433333333
$ws.Range("H12").Value = 56.79070299999999
$ws.Range("I12").Value = 0.04005030430848061
$ws.Range("J12").Value = 0.04005030430848062
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 9.844315000675778
$ws.Range("R12").Value = 88.598835006082
$ws.Range("S12").Value = 0.001942754217125311
$ws.Range("T12").Value = 0.001942754217125311
$ws.Range("G13").Value = 18.93023433333333
$ws.Range("H13").Value = 56.79070299999999
$ws.Range("I13").Value = 0.04005030430848061
$ws.Range("J13").Value = 0.04005030430848062
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 3.705732192468444
$ws.Range("R13").Value = 33.351589732216
$ws.Range("S13").Value = 0.000731318211979288
$ws.Range("T13").Value = 0.0007313182119792882
$ws.Range("G14").Value = 440.7369333333333
$ws.Range("H14").Value = 1322.2108
$ws.Range("I14").Value = 0.9324579922872165
$ws.Range("J14").Value = 0.9324579922872166
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 1216.233106870711
$ws.Range("R14").Value = 10946.0979618364
$ws.Range("S14").Value = 0.2400209661330719
$ws.Range("T14").Value = 0.2400209661330719
$ws.Range("G15").Value = 440.7369333333333
$ws.Range("H15").Value = 1322.2108
$ws.Range("I15").Value = 0.9324579922872165
$ws.Range("J15").Value = 0.9324579922872166
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 3193.239129283867
$ws.Range("R15").Value = 28739.15216355479
$ws.Range("S15").Value = 0.6301788173458414
$ws.Range("T15").Value = 0.6301788173458415
$ws.Range("G16").Value = 440.7369333333333
$ws.Range("H16").Value = 1322.2108
$ws.Range("I16").Value = 0.9324579922872165
$ws.Range("J16").Value = 0.9324579922872166
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 229.1970150905778
$ws.Range("R16").Value = 2062.7731358152
$ws.Range("S16").Value = 0.04523153389435294
$ws.Range("T16").Value = 0.04523153389435294
$ws.Range("G17").Value = 440.7369333333333
$ws.Range("H17").Value = 1322.2108
$ws.Range("I17").Value = 0.9324579922872165
$ws.Range("J17").Value = 0.9324579922872166
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 86.27748677084443
$ws.Range("R17").Value = 776.4973809375999
$ws.Range("S17").Value = 0.01702667491395034
$ws.Range("T17").Value = 0.01702667491395034
